$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update image names in column A (rows 2-5 change identity/order)
$ws.Range("A2").Value = "IMG_4402.JPG_4clusters.png"
$ws.Range("A3").Value = "IMG_4403.JPG_4clusters.png"
$ws.Range("A4").Value = "IMG_4399.JPG_4clusters.png"
$ws.Range("A5").Value = "IMG_4405.JPG_4clusters.png"

# Update Area_cm2 values in column B (rows 2-8)
$ws.Range("B2").Value = 2.749519586250239
$ws.Range("B3").Value = 1.726557732632442
$ws.Range("B4").Value = 2.878927310444024
$ws.Range("B5").Value = 5.092057187419351
$ws.Range("B6").Value = 4.964401783435124
$ws.Range("B7").Value = 6.01646226207128
$ws.Range("B8").Value = 4.445653220885297
